$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "29.287.95"
$cell.Style = "Normal"
$ws.Range("E2").Value = "  +0.45%  "
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "1.859.60"
$cell.Style = "Normal"
$ws.Range("E3").Value = "  +0.36%  "
$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = "1.000"
$cell.Style = "Normal"
$ws.Range("E4").Value = "  -0.01%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "0.7038"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  +2.22%  "
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "238.37"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  +0.24%  "
$ws.Range("E7").Value = "  +0.01%  "
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.08029"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  +4.12%  "
$ws.Range("E9").Value = "  -0.26%  "
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "23.57"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  +1.67%  "
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.08191"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  +0.52%  "
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "1.937.66"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  +4.27%  "
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "5.208"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  +0.02%  "
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "0.7091"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  -1.96%  "
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "89.75"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  +0.86%  "
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "29.344.96"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  +0.65%  "
$ws.Range("E17").Value = "  +1.88%  "
$ws.Range("E18").Value = "  +1.57%  "
$ws.Range("E19").Value = "  +0.98%  "
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "238.40"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  +1.59%  "
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "2.133.73"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  +1.55%  "
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "0.9987"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  -0.20%  "
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "1.000"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  +0.00%  "
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "7.478"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  -0.74%  "
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "163.07"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  +0.85%  "
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "8.898"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  -0.74%  "
$ws.Range("E27").Value = "  +1.48%  "
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "18.12"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  +0.37%  "
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "1.927"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  -1.73%  "
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "1.417"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  +1.39%  "
$ws.Range("E31").Value = "  -0.50%  "
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "4.381"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  -3.00%  "
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "4.029"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  +0.73%  "
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "0.05204"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  +0.36%  "
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "1.165"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  -0.91%  "
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "0.7190"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  +2.52%  "
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "1.006"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  -0.10%  "
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "2.691"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  +1.49%  "
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "0.01855"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  +0.06%  "
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "2.728"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  +1.97%  "
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "0.9443"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  +3.96%  "
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "1.153.73"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  +5.33%  "
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "6.002"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  +0.01%  "
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "0.4276"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  +0.17%  "
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "71.05"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  +0.83%  "
$ws.Range("E46").Value = "  -0.04%  "
$ws.Range("E47").Value = "  +0.14%  "
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "0.5300"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  -4.16%  "
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "1.765"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  +0.53%  "
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "2.029.04"
$cell.Style = "Normal"
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "9.177"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  +0.57%  "
